# Daily attendance processing - 2026-01-12 05:44:23
#
# The "Recorded By" column (G) lists the people/systems that recorded a
# session, separated by commas. For every row where it currently reads
# "dnasr281@gmail.com, System" the order of the two names should be
# swapped to "System, dnasr281@gmail.com". Rows that only contain a
# single name (e.g. just "System" or just "dnasr281@gmail.com") are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

# xlWhole = 2 (match the entire cell contents, not a substring)
$xlWhole = 2

$colG = $ws.Columns.Item(7)
$colG.Replace($oldValue, $newValue, $xlWhole) | Out-Null
